# Updated cryptos list on Sun May 21 07:26:31 UTC 2023 with GitHub Actions
# GitHub Actions bot refresh: pull the latest coinranking.com snapshot into
# the Price (D) / Volume(1h) (E) columns of the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.334.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "'1.834.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  +0.90%  "
$ws.Range("D5").Value = "'314.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D7").Value = "'0.4747"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.87%  "
$ws.Range("D8").Value = "'0.3692"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("D10").Value = "'0.8858"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("D11").Value = "'20.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("D12").Value = "'1.864.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").Value = "'0.07343"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.39%  "
$ws.Range("D14").Value = "'5.441"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'93.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.69%  "
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").Value = "'1.010"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").Value = "'0.000008798"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D20").Value = "'27.612.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.55%  "
$ws.Range("D21").Value = "'14.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("D22").Value = "'5.290"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'10.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("D24").Value = "'2.095.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("D25").Value = "'1.898"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "'152.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("D27").Value = "'18.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.94%  "
$ws.Range("D28").Value = "'2.146"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").Value = "'5.239"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "'117.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.54%  "
$ws.Range("D31").Value = "'0.08994"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("D32").Value = "'0.7511"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").Value = "'4.548"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("D35").Value = "'2.951"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.53%  "
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("D37").Value = "'1.104"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("D38").Value = "'0.05348"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("D39").Value = "'0.01957"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("D40").Value = "'2.974"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").Value = "'7.267"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.15%  "
$ws.Range("D42").Value = "'2.383"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.98%  "
$ws.Range("D43").Value = "'0.5314"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("D45").Value = "'8.489"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("D46").Value = "'0.4920"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("D47").Value = "'10.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("D48").Value = "'105.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.02%  "
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("D50").Value = "'1.673"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("D51").Value = "'0.06302"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.22%  "
